$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 487, shifting existing rows 487:612 down to 488:613
$ws.Rows.Item(487).Insert()

# Populate the newly inserted row 487 with the new data record
$ws.Cells.Item(487, 1).Value = 3
$ws.Cells.Item(487, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(487, 3).Value = "Coquimbo"
$ws.Cells.Item(487, 4).Value = 45135
$ws.Cells.Item(487, 5).Value = 5
$ws.Cells.Item(487, 6).Value = 100112040
$ws.Cells.Item(487, 7).Value = "Cilantro"
$ws.Cells.Item(487, 8).Value = "Sin especificar"
$ws.Cells.Item(487, 9).Value = "Primera"
$ws.Cells.Item(487, 10).Value = 210
$ws.Cells.Item(487, 11).Value = 3500
$ws.Cells.Item(487, 12).Value = 3800
$ws.Cells.Item(487, 13).Value = 3629
$ws.Cells.Item(487, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(487, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(487, 16).Value = 1210
$ws.Cells.Item(487, 17).Value = 3
$ws.Cells.Item(487, 18).Value = "Hortaliza"
